# Update the LR-pairs table (Tgfb1-Sdc2) with the recomputed NATMI values
# following Dr Hou advice: expand from 8 sending/target pairs (2 clusters)
# to the full 4x4 cluster cross product (ECs, FAPs, M2, sCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 16,20

$data[0,0] = "ECs"
$data[0,1] = "Tgfb1"
$data[0,2] = "Sdc2"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 96.320746
$data[0,7] = 288.962238
$data[0,8] = 0.3809824610908788
$data[0,9] = 0.3809824610908788
$data[0,10] = 2
$data[0,11] = 0.6666666666666666
$data[0,12] = 0.6646083333333334
$data[0,13] = 1.993825
$data[0,14] = 0.005515555660921567
$data[0,15] = 0.005515555660921567
$data[0,16] = 64.01557046448333
$data[0,17] = 576.1401341803501
$data[0,18] = 0.002101329969981627
$data[0,19] = 0.002101329969981627

$data[1,0] = "ECs"
$data[1,1] = "Tgfb1"
$data[1,2] = "Sdc2"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 96.320746
$data[1,7] = 288.962238
$data[1,8] = 0.3809824610908788
$data[1,9] = 0.3809824610908788
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 85.826024
$data[1,13] = 257.478072
$data[1,14] = 0.7122664414292983
$data[1,15] = 0.7122664414292983
$data[1,16] = 8266.826657893904
$data[1,17] = 74401.43992104514
$data[1,18] = 0.2713610218081763
$data[1,19] = 0.2713610218081763

$data[2,0] = "ECs"
$data[2,1] = "Tgfb1"
$data[2,2] = "Sdc2"
$data[2,3] = "M2"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 96.320746
$data[2,7] = 288.962238
$data[2,8] = 0.3809824610908788
$data[2,9] = 0.3809824610908788
$data[2,10] = 1
$data[2,11] = 0.3333333333333333
$data[2,12] = 0.009315666666666667
$data[2,13] = 0.027947
$data[2,14] = 0.00007731031261809587
$data[2,15] = 0.00007731031261809588
$data[2,16] = 0.8972919628206667
$data[2,17] = 8.075627665386
$data[2,18] = 0.00002945387316894738
$data[2,19] = 0.00002945387316894739

$data[3,0] = "ECs"
$data[3,1] = "Tgfb1"
$data[3,2] = "Sdc2"
$data[3,3] = "sCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 96.320746
$data[3,7] = 288.962238
$data[3,8] = 0.3809824610908788
$data[3,9] = 0.3809824610908788
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 33.99712866666667
$data[3,13] = 101.991386
$data[3,14] = 0.2821406925971621
$data[3,15] = 0.2821406925971621
$data[3,16] = 3274.628795031319
$data[3,17] = 29471.65915528187
$data[3,18] = 0.1074906554395519
$data[3,19] = 0.1074906554395519

$data[4,0] = "FAPs"
$data[4,1] = "Tgfb1"
$data[4,2] = "Sdc2"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 18.46467533333333
$data[4,7] = 55.394026
$data[4,8] = 0.07303429161291354
$data[4,9] = 0.07303429161291354
$data[4,10] = 2
$data[4,11] = 0.6666666666666666
$data[4,12] = 0.6646083333333334
$data[4,13] = 1.993825
$data[4,14] = 0.005515555660921567
$data[4,15] = 0.005515555660921567
$data[4,16] = 12.27177709882778
$data[4,17] = 110.44599388945
$data[4,18] = 0.0004028247005470018
$data[4,19] = 0.0004028247005470018

$data[5,0] = "FAPs"
$data[5,1] = "Tgfb1"
$data[5,2] = "Sdc2"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 18.46467533333333
$data[5,7] = 55.394026
$data[5,8] = 0.07303429161291354
$data[5,9] = 0.07303429161291354
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 85.826024
$data[5,13] = 257.478072
$data[5,14] = 0.7122664414292983
$data[5,15] = 0.7122664414292983
$data[5,16] = 1584.749668310875
$data[5,17] = 14262.74701479787
$data[5,18] = 0.05201987498943957
$data[5,19] = 0.05201987498943957

$data[6,0] = "FAPs"
$data[6,1] = "Tgfb1"
$data[6,2] = "Sdc2"
$data[6,3] = "M2"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 18.46467533333333
$data[6,7] = 55.394026
$data[6,8] = 0.07303429161291354
$data[6,9] = 0.07303429161291354
$data[6,10] = 1
$data[6,11] = 0.3333333333333333
$data[6,12] = 0.009315666666666667
$data[6,13] = 0.027947
$data[6,14] = 0.00007731031261809587
$data[6,15] = 0.00007731031261809588
$data[6,16] = 0.1720107605135555
$data[6,17] = 1.548096844622
$data[6,18] = 0.000005646303916435523
$data[6,19] = 0.000005646303916435524

$data[7,0] = "FAPs"
$data[7,1] = "Tgfb1"
$data[7,2] = "Sdc2"
$data[7,3] = "sCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 18.46467533333333
$data[7,7] = 55.394026
$data[7,8] = 0.07303429161291354
$data[7,9] = 0.07303429161291354
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 33.99712866666667
$data[7,13] = 101.991386
$data[7,14] = 0.2821406925971621
$data[7,15] = 0.2821406925971621
$data[7,16] = 627.7459430955596
$data[7,17] = 5649.713487860036
$data[7,18] = 0.02060594561901053
$data[7,19] = 0.02060594561901053

$data[8,0] = "M2"
$data[8,1] = "Tgfb1"
$data[8,2] = "Sdc2"
$data[8,3] = "ECs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 122.909391
$data[8,7] = 368.728173
$data[8,8] = 0.4861499128584522
$data[8,9] = 0.4861499128584522
$data[8,10] = 2
$data[8,11] = 0.6666666666666666
$data[8,12] = 0.6646083333333334
$data[8,13] = 1.993825
$data[8,14] = 0.005515555660921567
$data[8,15] = 0.005515555660921567
$data[8,16] = 81.68660550352499
$data[8,17] = 735.1794495317249
$data[8,18] = 0.002681386903922962
$data[8,19] = 0.002681386903922962

$data[9,0] = "M2"
$data[9,1] = "Tgfb1"
$data[9,2] = "Sdc2"
$data[9,3] = "FAPs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 122.909391
$data[9,7] = 368.728173
$data[9,8] = 0.4861499128584522
$data[9,9] = 0.4861499128584522
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 85.826024
$data[9,13] = 257.478072
$data[9,14] = 0.7122664414292983
$data[9,15] = 0.7122664414292983
$data[9,16] = 10548.82434179138
$data[9,17] = 94939.41907612245
$data[9,18] = 0.3462682684328532
$data[9,19] = 0.3462682684328532

$data[10,0] = "M2"
$data[10,1] = "Tgfb1"
$data[10,2] = "Sdc2"
$data[10,3] = "M2"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 122.909391
$data[10,7] = 368.728173
$data[10,8] = 0.4861499128584522
$data[10,9] = 0.4861499128584522
$data[10,10] = 1
$data[10,11] = 0.3333333333333333
$data[10,12] = 0.009315666666666667
$data[10,13] = 0.027947
$data[10,14] = 0.00007731031261809587
$data[10,15] = 0.00007731031261809588
$data[10,16] = 1.144982916759
$data[10,17] = 10.304846250831
$data[10,18] = 0.000037584401742347
$data[10,19] = 0.00003758440174234701

$data[11,0] = "M2"
$data[11,1] = "Tgfb1"
$data[11,2] = "Sdc2"
$data[11,3] = "sCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 122.909391
$data[11,7] = 368.728173
$data[11,8] = 0.4861499128584522
$data[11,9] = 0.4861499128584522
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 33.99712866666667
$data[11,13] = 101.991386
$data[11,14] = 0.2821406925971621
$data[11,15] = 0.2821406925971621
$data[11,16] = 4178.566380168641
$data[11,17] = 37607.09742151778
$data[11,18] = 0.1371626731199337
$data[11,19] = 0.1371626731199337

$data[12,0] = "sCs"
$data[12,1] = "Tgfb1"
$data[12,2] = "Sdc2"
$data[12,3] = "ECs"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 15.127183
$data[12,7] = 45.381549
$data[12,8] = 0.05983333443775553
$data[12,9] = 0.05983333443775553
$data[12,10] = 2
$data[12,11] = 0.6666666666666666
$data[12,12] = 0.6646083333333334
$data[12,13] = 1.993825
$data[12,14] = 0.005515555660921567
$data[12,15] = 0.005515555660921567
$data[12,16] = 10.05365188165833
$data[12,17] = 90.482866934925
$data[12,18] = 0.0003300140864699758
$data[12,19] = 0.0003300140864699758

$data[13,0] = "sCs"
$data[13,1] = "Tgfb1"
$data[13,2] = "Sdc2"
$data[13,3] = "FAPs"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 15.127183
$data[13,7] = 45.381549
$data[13,8] = 0.05983333443775553
$data[13,9] = 0.05983333443775553
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 85.826024
$data[13,13] = 257.478072
$data[13,14] = 0.7122664414292983
$data[13,15] = 0.7122664414292983
$data[13,16] = 1298.305971210392
$data[13,17] = 11684.75374089353
$data[13,18] = 0.04261727619882921
$data[13,19] = 0.04261727619882921

$data[14,0] = "sCs"
$data[14,1] = "Tgfb1"
$data[14,2] = "Sdc2"
$data[14,3] = "M2"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 15.127183
$data[14,7] = 45.381549
$data[14,8] = 0.05983333443775553
$data[14,9] = 0.05983333443775553
$data[14,10] = 1
$data[14,11] = 0.3333333333333333
$data[14,12] = 0.009315666666666667
$data[14,13] = 0.027947
$data[14,14] = 0.00007731031261809587
$data[14,15] = 0.00007731031261809588
$data[14,16] = 0.1409197944336667
$data[14,17] = 1.268278149903
$data[14,18] = 0.000004625733790365961
$data[14,19] = 0.000004625733790365962

$data[15,0] = "sCs"
$data[15,1] = "Tgfb1"
$data[15,2] = "Sdc2"
$data[15,3] = "sCs"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 15.127183
$data[15,7] = 45.381549
$data[15,8] = 0.05983333443775553
$data[15,9] = 0.05983333443775553
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 33.99712866666667
$data[15,13] = 101.991386
$data[15,14] = 0.2821406925971621
$data[15,15] = 0.2821406925971621
$data[15,16] = 514.2807868152128
$data[15,17] = 4628.527081336914
$data[15,18] = 0.01688141841866597
$data[15,19] = 0.01688141841866597

$ws.Range("A2:T17").Value = $data
Write-Output "Updated rows 2:17 (16 data rows) of Sheet1"
